# The document uses a distinct "first page" header/footer (titlePg) as well
# as the default/primary header/footer, giving four logo pictures in total:
#   - Sections(1).Headers(1) -> default header    (BTec logo,    docPr id="3")
#   - Sections(1).Headers(2) -> first-page header  (BTec logo,    docPr id="1")
#   - Sections(1).Footers(1) -> default footer     (Pearson logo, docPr id="4")
#   - Sections(1).Footers(2) -> first-page footer  (Pearson logo, docPr id="2")
#
# Each of those inline pictures carries a docPr "name" that was mislabelled:
#   - the BTEC logo was named "image2.jpg" and should be "image1.jpg"
#   - the Pearson logo was named "image1.png" and should be "image2.png"

$d = $word.ActiveDocument

$headerTarget = "image1.jpg"
$footerTarget = "image2.png"

$h1 = $d.Sections(1).Headers(1).Range.InlineShapes(1)
$h1.Name = $headerTarget

$h2 = $d.Sections(1).Headers(2).Range.InlineShapes(1)
$h2.Name = $headerTarget

# Footers hold several paragraphs of text before the logo paragraph, so the
# picture is reached through its own paragraph range rather than indexing
# InlineShapes directly off the footer's full Range.
$f1 = $d.Sections(1).Footers(1)
$f1Pic = $f1.Range.Paragraphs($f1.Range.Paragraphs.Count).Range.InlineShapes(1)
$f1Pic.Name = $footerTarget

$f2 = $d.Sections(1).Footers(2)
$f2Pic = $f2.Range.Paragraphs($f2.Range.Paragraphs.Count).Range.InlineShapes(1)
$f2Pic.Name = $footerTarget
